$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook lists weekly price observations for "Zapallo italiano" (Primera /
# Segunda quality) starting at row 81. A new pair of observations (for date
# 44469) needs to be inserted right after the most recent existing entries
# (rows 79-80, date 44446) and before the historical entries that currently
# start at row 81 (date 44299). Inserting 2 rows at row 81 shifts every
# existing row 81-209 down to 83-211, matching the target dimension A1:R211.
$ws.Range("A81:A82").EntireRow.Insert()

# ---- Fill new row 81 (Primera) ----
$ws.Cells.Item(81, 1).Value = 1
$ws.Cells.Item(81, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(81, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(81, 4).Value = 44469
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(83, 4).NumberFormat
$ws.Cells.Item(81, 5).Value = 15
$ws.Cells.Item(81, 6).Value = 100112032
$ws.Cells.Item(81, 7).Value = "Zapallo italiano"
$ws.Cells.Item(81, 8).Value = "Huracán"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 120
$ws.Cells.Item(81, 11).Value = 9000
$ws.Cells.Item(81, 12).Value = 10000
$ws.Cells.Item(81, 13).Value = 9500
$ws.Cells.Item(81, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(81, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(81, 16).Value = 136
$ws.Cells.Item(81, 17).Value = 70
$ws.Cells.Item(81, 18).Value = "Hortaliza"

# ---- Fill new row 82 (Segunda) ----
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value = 44469
$ws.Cells.Item(82, 4).NumberFormat = $ws.Cells.Item(83, 4).NumberFormat
$ws.Cells.Item(82, 5).Value = 15
$ws.Cells.Item(82, 6).Value = 100112032
$ws.Cells.Item(82, 7).Value = "Zapallo italiano"
$ws.Cells.Item(82, 8).Value = "Huracán"
$ws.Cells.Item(82, 9).Value = "Segunda"
$ws.Cells.Item(82, 10).Value = 120
$ws.Cells.Item(82, 11).Value = 9000
$ws.Cells.Item(82, 12).Value = 10000
$ws.Cells.Item(82, 13).Value = 9500
$ws.Cells.Item(82, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 95
$ws.Cells.Item(82, 17).Value = 100
$ws.Cells.Item(82, 18).Value = "Hortaliza"
